$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Drop the stray "_GoBack" bookmark that used to sit in the middle of the
#    "At the bottom of the screen ... This screen allows ..." paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Append a status-update note at the very end of the document: a blank
#    separator paragraph followed by a new paragraph of text, both justified
#    like the rest of the write-up, and finish by re-planting "_GoBack" at
#    the new end of the document (where Word leaves it after the last edit).
# ---------------------------------------------------------------------------
$noteText = "Note: The navigation and entering text has been implemented. " + `
    "I am currently looking up how to open up the camera and library on pressing the corresponding buttons. I have not been able to " + `
    "implement that yet." + `
    " I am also still playing with Realm to understand it better and will most probably be sticking to that for data persistence. "

$lastPara = $d.Paragraphs.Last
$endOfDoc = $lastPara.Range
$endOfDoc.Collapse(0)

# Type the new paragraph in one shot. A one-character sentinel ("Z") is
# tacked on the end and removed afterwards - it keeps the freshly-typed text
# from ever being the literal last character of the document while we are
# still positioning the bookmark, which sidesteps an end-of-story addressing
# quirk when re-resolving a Range right at Content.End.
$endOfDoc.InsertAfter("`r" + $noteText + "Z")

# Paragraph formatting (centered/both-justified like its neighbours) is
# inherited automatically from the paragraph mark it split off from, so the
# new paragraph already has <w:jc w:val="both"/>.

# Now slot a blank paragraph in right before the text we just typed, so the
# final layout is: ... Realm paragraph / <blank paragraph> / new note
# paragraph. Targeting the *start* of the new paragraph (rather than the end
# of the previous one) keeps this insertion from fusing into the Realm
# paragraph's run.
$newLastPara = $d.Paragraphs.Last
$splitPoint = $newLastPara.Range.Start
$d.Range($splitPoint, $splitPoint).InsertAfter("`r")

# Re-create the "_GoBack" bookmark as a zero-length bookmark right after the
# new note text (immediately before the trailing sentinel, which sits
# immediately before the final paragraph mark).
$bookmarkPos = $d.Content.End - 2
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))

# Remove the temporary sentinel character now that the bookmark is anchored.
$sentinelRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$sentinelRange.Delete()
